$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 800
$ws.Range("I62").Value = 800
$ws.Range("K62").Value = 800
$ws.Range("M62").Value = -176

$ws.Range("H65").Value = 800
$ws.Range("I65").Value = 800
$ws.Range("K65").Value = 4000
$ws.Range("M65").Value = -880

$ws.Range("H87").Value = 109899.5
$ws.Range("J87").Value = 109899.5
$ws.Range("L87").Value = 109899.5
$ws.Range("N87").Value = -112395.5

$ws.Range("H90").Value = 109899.5
$ws.Range("J90").Value = 109899.5
$ws.Range("L90").Value = 329698.5
$ws.Range("N90").Value = -342178.5

$ws.Range("H98").Value = 22699.5
$ws.Range("I98").Value = 19570.715
$ws.Range("K98").Value = 19570.715
$ws.Range("M98").Value = -18072.715

$ws.Range("H113").Value = 2002.5
$ws.Range("I113").Value = 2002.5
$ws.Range("K113").Value = 2002.5
$ws.Range("M113").Value = 1251.5

$ws.Range("H122").Value = 22699.5
$ws.Range("I122").Value = 19570.715
$ws.Range("K122").Value = 58712.145
$ws.Range("M122").Value = -56262.145

$ws.Range("H138").Value = 1691.875
$ws.Range("I138").Value = 922.8333
$ws.Range("J138").Value = 3999
$ws.Range("K138").Value = 2768.4999
$ws.Range("L138").Value = 11997
$ws.Range("M138").Value = 2371.5001
$ws.Range("N138").Value = -22277

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 759.2
$ws.Range("I2").Value = 732.44446
$ws.Range("K2").Value = 732.44446
$ws.Range("M2").Value = -619.44446

$ws.Range("H5").Value = 35.5
$ws.Range("I5").Value = 35.5
$ws.Range("K5").Value = 35.5
$ws.Range("M5").Value = 76.5

$ws.Range("H92").Value = 47516.5
$ws.Range("J92").Value = 47516.5
$ws.Range("L92").Value = 47516.5
$ws.Range("N92").Value = -52508.5

$ws.Range("H116").Value = 759.2
$ws.Range("I116").Value = 732.44446
$ws.Range("K116").Value = 732.44446
$ws.Range("M116").Value = 1561.55554

$ws.Range("H132").Value = 7037.143
$ws.Range("I132").Value = 7037.143
$ws.Range("K132").Value = 21111.429
$ws.Range("M132").Value = -18581.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 759.2
$ws.Range("I3").Value = 732.44446
$ws.Range("K3").Value = 732.44446
$ws.Range("M3").Value = -618.44446

$ws.Range("H4").Value = 35.5
$ws.Range("I4").Value = 35.5
$ws.Range("K4").Value = 35.5
$ws.Range("M4").Value = 79.5

$ws.Range("H86").Value = 3005.1428
$ws.Range("I86").Value = 2790.6
$ws.Range("J86").Value = 3541.5
$ws.Range("K86").Value = 2790.6
$ws.Range("L86").Value = 3541.5
$ws.Range("M86").Value = -1667.6
$ws.Range("N86").Value = -5787.5

$ws.Range("H89").Value = 3005.1428
$ws.Range("I89").Value = 2790.6
$ws.Range("J89").Value = 3541.5
$ws.Range("K89").Value = 13953
$ws.Range("L89").Value = 17707.5
$ws.Range("M89").Value = -8337
$ws.Range("N89").Value = -28939.5

$ws.Range("H92").Value = 98400.336
$ws.Range("J92").Value = 98400.336
$ws.Range("L92").Value = 98400.336
$ws.Range("N92").Value = -103392.336

$ws.Range("H99").Value = 3839.6
$ws.Range("I99").Value = 3839.6
$ws.Range("K99").Value = 3839.6
$ws.Range("M99").Value = -2341.6

$ws.Range("H105").Value = 1342.7142
$ws.Range("I105").Value = 1316.5
$ws.Range("K105").Value = 1316.5
$ws.Range("M105").Value = 430.5

$ws.Range("H134").Value = 1000
$ws.Range("I134").Value = 1000
$ws.Range("K134").Value = 3000
$ws.Range("M134").Value = -465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1200
$ws.Range("I16").Value = 1200
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1200
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -913
$ws.Range("N16").ClearContents()

$ws.Range("H107").Value = 175.3158
$ws.Range("I107").Value = 198.85715
$ws.Range("K107").Value = 198.85715
$ws.Range("M107").Value = 1721.14285

$ws.Range("H113").Value = 1200
$ws.Range("I113").Value = 1200
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1200
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 970
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 833
$ws.Range("I132").Value = 833
$ws.Range("K132").Value = 2499
$ws.Range("M132").Value = 31

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1979
$ws.Range("I117").Value = 575
$ws.Range("J117").Value = 2681
$ws.Range("K117").Value = 1725
$ws.Range("L117").Value = 8043
$ws.Range("M117").Value = 1717
$ws.Range("N117").Value = -14927

$ws.Range("H119").Value = 1107.25
$ws.Range("I119").Value = 1107.25
$ws.Range("K119").Value = 3321.75
$ws.Range("M119").Value = 1516.25

$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()

$ws.Range("H124").Value = 500
$ws.Range("I124").Value = 500
$ws.Range("K124").Value = 1500
$ws.Range("M124").Value = 3410

$ws.Range("H125").Value = 1000
$ws.Range("I125").Value = 1000
$ws.Range("K125").Value = 3000
$ws.Range("M125").Value = 1920

$ws.Range("H131").Value = 5750
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 5750
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 17250
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -27330

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 24667.705
$ws.Range("I102").Value = 24667.705
$ws.Range("K102").Value = 24667.705
$ws.Range("M102").Value = -23045.705

$ws.Range("H107").Value = 1709.1
$ws.Range("I107").Value = 929.8
$ws.Range("J107").Value = 2488.4
$ws.Range("K107").Value = 929.8
$ws.Range("L107").Value = 2488.4
$ws.Range("M107").Value = 990.2
$ws.Range("N107").Value = -6328.4

$ws.Range("H123").Value = 89497
$ws.Range("J123").Value = 89497
$ws.Range("L123").Value = 89497
$ws.Range("N123").Value = -94397

$ws.Range("H132").Value = 3200
$ws.Range("I132").Value = 3200
$ws.Range("K132").Value = 9600
$ws.Range("M132").Value = -7070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2699.4546
$ws.Range("I22").Value = 1099.3334
$ws.Range("J22").Value = 3299.5
$ws.Range("K22").Value = 1099.3334
$ws.Range("L22").Value = 3299.5
$ws.Range("M22").Value = -804.3334
$ws.Range("N22").Value = -3889.5

$ws.Range("H27").Value = 2699.4546
$ws.Range("I27").Value = 1099.3334
$ws.Range("J27").Value = 3299.5
$ws.Range("K27").Value = 1099.3334
$ws.Range("L27").Value = 3299.5
$ws.Range("M27").Value = -992.3334
$ws.Range("N27").Value = -3513.5

$ws.Range("H46").Value = 3283.7144
$ws.Range("I46").Value = 944
$ws.Range("J46").Value = 4219.6
$ws.Range("K46").Value = 944
$ws.Range("L46").Value = 4219.6
$ws.Range("M46").Value = -756
$ws.Range("N46").Value = -4595.6

$ws.Range("H61").Value = 3607.1428
$ws.Range("I61").Value = 4025
$ws.Range("J61").Value = 1100
$ws.Range("K61").Value = 4025
$ws.Range("L61").Value = 1100
$ws.Range("M61").Value = -3823
$ws.Range("N61").Value = -1504

$ws.Range("H113").Value = 3607.1428
$ws.Range("I113").Value = 4025
$ws.Range("J113").Value = 1100
$ws.Range("K113").Value = 4025
$ws.Range("L113").Value = 1100
$ws.Range("M113").Value = -1855
$ws.Range("N113").Value = -5440

$ws.Range("H132").Value = 2749.5
$ws.Range("J132").Value = 3500
$ws.Range("L132").Value = 10500
$ws.Range("N132").Value = -15560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 9210.25
$ws.Range("I64").Value = 5263
$ws.Range("K64").Value = 5263
$ws.Range("M64").Value = -5015

$ws.Range("H67").Value = 9210.25
$ws.Range("I67").Value = 5263
$ws.Range("K67").Value = 5263
$ws.Range("M67").Value = -4405

$ws.Range("H107").Value = 403
$ws.Range("I107").Value = 353.75
$ws.Range("K107").Value = 1061.25
$ws.Range("M107").Value = 858.75

Write-Output "Applied all profit updates"